$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("L2").Value = 5930
$ws.Range("L3").Value = 6454
$ws.Range("L4").Value = 1588
$ws.Range("L6").Value = 5311
$ws.Range("L7").Value = 19666

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("L2").Value = 173
$ws.Range("L5").Value = 72
$ws.Range("L8").Value = 1294
$ws.Range("L10").Value = 131
$ws.Range("L14").Value = 98
$ws.Range("L15").Value = 160
$ws.Range("L19").Value = 536
$ws.Range("L25").Value = 117
$ws.Range("L29").Value = 1106
$ws.Range("L34").Value = 111
$ws.Range("L36").Value = 247
$ws.Range("L37").Value = 746
$ws.Range("L42").Value = 633
$ws.Range("L47").Value = 137
$ws.Range("L49").Value = 105
$ws.Range("L52").Value = 413
$ws.Range("L57").Value = 68
$ws.Range("L63").Value = 57
$ws.Range("L65").Value = 384
$ws.Range("L67").Value = 679
$ws.Range("L75").Value = 69
$ws.Range("L76").Value = 300
$ws.Range("L77").Value = 133
$ws.Range("L83").Value = 429
$ws.Range("L84").Value = 188
$ws.Range("L85").Value = 976
$ws.Range("L93").Value = 99
$ws.Range("L94").Value = 244
$ws.Range("L99").Value = 341
$ws.Range("L101").Value = 19666

$ws = $wb.Worksheets.Item("Bridgeport")
$ws.Range("L2").Value = 40
$ws.Range("L7").Value = 98

$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("L3").Value = 202
$ws.Range("L4").Value = 45

$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("L3").Value = 404
$ws.Range("L6").Value = 204
$ws.Range("L7").Value = 976

$ws = $wb.Worksheets.Item("Little Village")
$ws.Range("L3").Value = 129
$ws.Range("L4").Value = 27
$ws.Range("L6").Value = 117
$ws.Range("L7").Value = 413

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("L2").Value = 388
$ws.Range("L3").Value = 458
$ws.Range("L7").Value = 1294

$ws = $wb.Worksheets.Item("South Chicago")
$ws.Range("L2").Value = 137
$ws.Range("L3").Value = 173
$ws.Range("L7").Value = 429

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("L2").Value = 224
$ws.Range("L7").Value = 746

$ws = $wb.Worksheets.Item("New City")
$ws.Range("L4").Value = 21
$ws.Range("L7").Value = 384

$ws = $wb.Worksheets.Item("Woodlawn")
$ws.Range("L2").Value = 97
$ws.Range("L7").Value = 341

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("L3").Value = 264
$ws.Range("L6").Value = 159
$ws.Range("L7").Value = 679

$ws = $wb.Worksheets.Item("South Deering")
$ws.Range("L6").Value = 56
$ws.Range("L7").Value = 188

$ws = $wb.Worksheets.Item("Lincoln Park")
$ws.Range("L4").Value = 14
$ws.Range("L7").Value = 105

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("L2").Value = 326
$ws.Range("L3").Value = 430
$ws.Range("L6").Value = 271
$ws.Range("L7").Value = 1106

$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("L2").Value = 190
$ws.Range("L4").Value = 26
$ws.Range("L7").Value = 536

$ws = $wb.Worksheets.Item("River North")
$ws.Range("L3").Value = 59
$ws.Range("L7").Value = 300

$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("L2").Value = 172
$ws.Range("L6").Value = 176
$ws.Range("L7").Value = 633

$ws = $wb.Worksheets.Item("Avondale")
$ws.Range("L3").Value = 38
$ws.Range("L7").Value = 131

$ws = $wb.Worksheets.Item("Grand Boulevard")
$ws.Range("L2").Value = 87
$ws.Range("L7").Value = 247

$ws = $wb.Worksheets.Item("West Lawn")
$ws.Range("L3").Value = 26
$ws.Range("L7").Value = 99

$ws = $wb.Worksheets.Item("Garfield Ridge")
$ws.Range("L2").Value = 39
$ws.Range("L7").Value = 111

$ws = $wb.Worksheets.Item("West Loop")
$ws.Range("L6").Value = 92
$ws.Range("L7").Value = 244

$ws = $wb.Worksheets.Item("East Side")
$ws.Range("L3").Value = 54
$ws.Range("L7").Value = 117

$ws = $wb.Worksheets.Item("Kenwood")
$ws.Range("L3").Value = 46
$ws.Range("L7").Value = 137

$ws = $wb.Worksheets.Item("Brighton Park")
$ws.Range("L3").Value = 52
$ws.Range("L7").Value = 160

$ws = $wb.Worksheets.Item("Portage Park")
$ws.Range("L2").Value = 57
$ws.Range("L4").Value = 14

$ws = $wb.Worksheets.Item("Albany Park")
$ws.Range("L3").Value = 55
$ws.Range("L4").Value = 14
$ws.Range("L7").Value = 173

$ws = $wb.Worksheets.Item("Armour Square")
$ws.Range("L3").Value = 20
$ws.Range("L7").Value = 72

$ws = $wb.Worksheets.Item("Pullman")
$ws.Range("L2").Value = 30
$ws.Range("L7").Value = 69

$ws = $wb.Worksheets.Item("Mckinley Park")
$ws.Range("L2").Value = 22
$ws.Range("L7").Value = 68

$ws = $wb.Worksheets.Item("Riverdale")
$ws.Range("L2").Value = 51
$ws.Range("L7").Value = 133
